$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on D2:E51 so numeric-looking strings (e.g. "355.82") are kept as text
# instead of being parsed into floating point numbers, matching the inlineStr text cells
# in the original workbook.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "51.700.59"
$ws.Range("E2").Value = "  -0.36%  "

# Row 3
$ws.Range("D3").Value = "2.776.04"
$ws.Range("E3").Value = "  -1.27%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "355.82"
$ws.Range("E5").Value = "  +0.55%  "

# Row 6
$ws.Range("D6").Value = "108.70"
$ws.Range("E6").Value = "  -2.36%  "

# Row 7
$ws.Range("D7").Value = "0.554"
$ws.Range("E7").Value = "  -1.92%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").Value = "0.584"

# Row 10
$ws.Range("D10").Value = "39.57"
$ws.Range("E10").Value = "  -2.11%  "

# Row 11
$ws.Range("E11").Value = "  +2.73%  "

# Row 12
$ws.Range("D12").Value = "0.0844"
$ws.Range("E12").Value = "  -1.14%  "

# Row 13
$ws.Range("D13").Value = "19.38"
$ws.Range("E13").Value = "  -2.34%  "

# Row 14
$ws.Range("D14").Value = "7.58"
$ws.Range("E14").Value = "  -2.00%  "

# Row 15
$ws.Range("D15").Value = "3.211.83"
$ws.Range("E15").Value = "  -1.34%  "

# Row 16
$ws.Range("D16").Value = "2.780.16"
$ws.Range("E16").Value = "  -1.33%  "

# Row 17
$ws.Range("D17").Value = "0.929"
$ws.Range("E17").Value = "  +0.54%  "

# Row 18
$ws.Range("D18").Value = "51.670.38"
$ws.Range("E18").Value = "  -0.08%  "

# Row 19
$ws.Range("D19").Value = "7.40"
$ws.Range("E19").Value = "  -1.28%  "

# Row 20
$ws.Range("D20").Value = "3.08"
$ws.Range("E20").Value = "  -1.21%  "

# Row 21
$ws.Range("D21").Value = "13.06"
$ws.Range("E21").Value = "  -1.94%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0967"
$ws.Range("E22").Value = "  -2.34%  "

# Row 23
$ws.Range("D23").Value = "70.03"
$ws.Range("E23").Value = "  -0.50%  "

# Row 24
$ws.Range("D24").Value = "267.85"
$ws.Range("E24").Value = "  +0.14%  "

# Row 25
$ws.Range("E25").Value = "  -2.19%  "

# Row 26
$ws.Range("D26").Value = "26.33"
$ws.Range("E26").Value = "  -1.88%  "

# Row 27
$ws.Range("E27").Value = "  -0.05%  "

# Row 28
$ws.Range("E28").Value = "  +17.40%  "

# Row 29
$ws.Range("E29").Value = "  +1.61%  "

# Row 30
$ws.Range("D30").Value = "10.17"
$ws.Range("E30").Value = "  -1.01%  "

# Row 31
$ws.Range("D31").Value = "6.19"
$ws.Range("E31").Value = "  +4.90%  "

# Row 32
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "51.59"
$ws.Range("E32").Value = "  -1.82%  "

# Row 33
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "34.62"
$ws.Range("E33").Value = "  +0.89%  "

# Row 34
$ws.Range("E34").Value = "  -7.81%  "

# Row 35
$ws.Range("D35").Value = "0.0837"
$ws.Range("E35").Value = "  -1.18%  "

# Row 36
$ws.Range("E36").Value = "  -7.38%  "

# Row 37
$ws.Range("E37").Value = "  +0.06%  "

# Row 38
$ws.Range("D38").Value = "18.55"
$ws.Range("E38").Value = "  +1.57%  "

# Row 39
$ws.Range("D39").Value = "3.11"
$ws.Range("E39").Value = "  -3.30%  "

# Row 40
$ws.Range("E40").Value = "  -3.63%  "

# Row 41
$ws.Range("D41").Value = "2.54"
$ws.Range("E41").Value = "  +2.55%  "

# Row 42
$ws.Range("E42").Value = "  -2.54%  "

# Row 43
$ws.Range("E43").Value = "  -2.47%  "

# Row 44
$ws.Range("D44").Value = "118.89"
$ws.Range("E44").Value = "  -5.81%  "

# Row 45
$ws.Range("D45").Value = "21.73"
$ws.Range("E45").Value = "  -5.56%  "

# Row 46
$ws.Range("D46").Value = "2.079.37"
$ws.Range("E46").Value = "  +0.09%  "

# Row 47
$ws.Range("D47").Value = "3.25"
$ws.Range("E47").Value = "  -2.17%  "

# Row 48
$ws.Range("E48").Value = "  +1.04%  "

# Row 49
$ws.Range("D49").Value = "0.937"
$ws.Range("E49").Value = "  -4.00%  "

# Row 50
$ws.Range("E50").Value = "  -5.82%  "

# Row 51
$ws.Range("D51").Value = "8.53"
$ws.Range("E51").Value = "  -6.05%  "

# Restore default (Normal) style so no stray number-format styling is left on the cells
$ws.Range("D2:E51").Style = "Normal"